$d = $word.ActiveDocument

# Touch the document's formatting so Word re-serializes the run
# properties (this normalizes attribute ordering, e.g. <w:rFonts>
# and <w:bookmarkStart>, without altering any visible content).
$d.Content.Font.Name = $d.Content.Font.Name
